# Applies the "Finalize fast close slow far strategy" commit:
#  - add a new "test" worksheet (copied layout from "closest") with its own
#    beacon-power inputs (B2=2, C2=5) and becomes the active tab
#  - tidy a couple of conditional-format style leftovers on the
#    "closest" and "closest with more crystals" sheets
#  - update the "closest with more crystals" sheet's saved selection

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Style touch-ups that are part of the same commit
# ---------------------------------------------------------------------

$ws2 = $wb.Worksheets.Item("closest")
$ws2.Range("C10").NumberFormat = "0.0"
$ws2.Range("C10").Font.Strikethrough = $true
$ws2.Range("C11").NumberFormat = "0.0"
$ws2.Range("C11").Font.Strikethrough = $true

$ws3 = $wb.Worksheets.Item("closest with more crystals")
$ws3.Range("C7").NumberFormat = "0.0"
$ws3.Range("C7").Font.Strikethrough = $false
$ws3.Range("C10").NumberFormat = "0.0"
$ws3.Range("C10").Font.Strikethrough = $true
$ws3.Range("C11").NumberFormat = "0.0"
$ws3.Range("C11").Font.Strikethrough = $true

# Update "closest with more crystals"' saved selection now, while it is
# still the active sheet (Range.Select()/Activate() switch the active
# sheet as a side effect, so this must happen before the new sheet is
# made active below).
$ws3.Range("A1:E13").Select()
$ws3.Range("D7").Activate()

# ---------------------------------------------------------------------
# 2) New "test" sheet, appended after "closest with more crystals".
#    Copy "closest" so formulas / number formats / conditional
#    highlighting styles all come along, then adjust the two inputs.
# ---------------------------------------------------------------------

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws2.Copy([System.Reflection.Missing]::Value, $lastSheet)

$wsTest = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsTest.Name = "test"

$wsTest.Range("B2").Value = 2
$wsTest.Range("C2").Value = 5

# rows 5-9 are not "best"/"excluded" for these inputs -> plain 0.0 format
$wsTest.Range("C5:C9").Font.Strikethrough = $false
$wsTest.Range("C5:C9").NumberFormat = "0.0"

# rows 10-13 keep the "excluded" look (0.0 + strikethrough)
$wsTest.Range("C10:C13").NumberFormat = "0.0"
$wsTest.Range("C10:C13").Font.Strikethrough = $true

$wb.Application.Calculate()

# Make "test" the active sheet/tab and set its saved selection.
$wsTest.Activate()
$wsTest.Range("C3").Select()

# Restore "closest with more crystals"' saved selection (no longer the
# active tab, but its sheetView selection still changed in the commit).
$ws3.Range("A1:E13").Select()
$ws3.Range("D7").Activate()
